$d = $word.ActiveDocument
$d.Content.Find.Execute("Docker JavaScript", $false, $false, $false, $false, $false, $true, 1, $false, "Docker, JavaScript", 2)
